$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.862.04'
$ws.Range('E2').Value = '  -0.26%  '

$ws.Range('D3').Value = '2.414.84'
$ws.Range('E3').Value = '  +0.10%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.75'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.82%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.35'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.23%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.530'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('E9').Value = '  +0.24%  '

$ws.Range('E10').Value = '  -1.23%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.20'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.63%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.349'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.08%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.69'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.31%  '

$ws.Range('E14').Value = '  -0.99%  '

$ws.Range('D15').Value = '2.851.73'
$ws.Range('E15').Value = '  +0.44%  '

$ws.Range('D16').Value = '61.846.51'
$ws.Range('E16').Value = '  -0.06%  '

$ws.Range('D17').Value = '2.410.79'
$ws.Range('E17').Value = '  +0.11%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.26'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.12%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '323.03'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.79'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.85%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.11'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.73%  '

$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.51'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.72'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.75'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.73%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '553.48'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.70%  '

$ws.Range('E27').Value = '  +0.48%  '

$ws.Range('E28').Value = '  -0.06%  '

$ws.Range('D29').Value = '0.0₃0929'
$ws.Range('E29').Value = '  -0.34%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.14'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.40%  '

$ws.Range('E31').Value = '  -4.90%  '

$ws.Range('E32').Value = '  -0.93%  '

$ws.Range('E33').Value = '  -0.64%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.49'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.53%  '

$ws.Range('E35').Value = '  +0.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.72'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.78%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.378'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.39%  '

$ws.Range('E38').Value = '  +1.92%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.40'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.87%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.47'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.79'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.46%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.991'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.82%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.23'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.44%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '146.70'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.76%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.62'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.66%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0524'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.03%  '

$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.591'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.48%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.70'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.59%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0919'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.45%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0226'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.95%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.55'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.69%  '
